$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.518.83"
$ws.Range("E2").Value = "  +1.80%  "
$ws.Range("D3").Value = "2.682.83"
$ws.Range("E3").Value = "  +2.24%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.53%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +0.35%  "
$ws.Range("D9").Value = "2.681.04"
$ws.Range("E9").Value = "  +2.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.171"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.05%  "
$ws.Range("E11").Value = "  +2.25%  "
$ws.Range("E12").Value = "  +2.77%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.06"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.67%  "
$ws.Range("D14").Value = "3.170.77"
$ws.Range("E14").Value = "  +3.01%  "
$ws.Range("E15").Value = "  +1.26%  "
$ws.Range("D16").Value = "72.376.74"
$ws.Range("E16").Value = "  +1.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.40"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.65%  "
$ws.Range("D18").Value = "2.690.72"
$ws.Range("E18").Value = "  +2.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.88%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.05"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.87%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "373.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.48%  "
$ws.Range("E22").Value = "  +1.20%  "
$ws.Range("E23").Value = "  +9.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "72.46"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.35%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("E26").Value = "  -2.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.93"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.97%  "
$ws.Range("D28").Value = "2.819.21"
$ws.Range("E28").Value = "  +2.17%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.09"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.75%  "
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "518.86"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.67%  "
$ws.Range("E33").Value = "  -0.93%  "
$ws.Range("E34").Value = "  -0.78%  "
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "164.63"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.60%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.56"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.94%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.13"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.58%  "
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("E40").Value = "  -2.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.109"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.42%  "
$ws.Range("E43").Value = "  -0.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.59"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.31%  "
$ws.Range("E45").Value = "  +1.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.24"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "154.37"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.09%  "
$ws.Range("E48").Value = "  +2.70%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.551"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.51%  "
$ws.Range("E50").Value = "  +2.43%  "
$ws.Range("E51").Value = "  +1.77%  "
